# 17.4.1.xlsx update: extend the year table from column S (2022) through
# column U (2024), i.e. add two more year columns (2023, 2024) with their
# data-series values, extend row3's bottom border across the new columns,
# widen the new columns to match the rest of the data columns, bump the
# row-5 height so the wrapped header text fits, and move the selection
# back to the default cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (the thin border-only row above the header) -------------------
# Clone S3's formatting (border only, no value) into the two new cells.
$ws.Range("S3").Copy($ws.Range("T3"))
$ws.Range("S3").Copy($ws.Range("U3"))

# --- Row 4 (the year header row) -----------------------------------------
$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("S4").Copy($ws.Range("U4"))
$ws.Range("T4").Value = 2023
$ws.Range("U4").Value = 2024

# --- Row 5 (the data row) --------------------------------------------------
$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("S5").Copy($ws.Range("U5"))
$ws.Range("T5").Value = 10.8
$ws.Range("U5").Value = 6.5

# Row 5 grows a little taller so the (wrapped) header text still fits next
# to the two extra columns.
$ws.Rows.Item(5).RowHeight = 41.25

# --- Column widths ----------------------------------------------------------
# Columns D..U (4..21) all share the same custom width used by the rest of
# the numeric columns.
$ws.Range("D1:U1").ColumnWidth = 7.83

# --- Selection / view -------------------------------------------------------
# Reset the active cell back to the top-left of the sheet (the workbook had
# been left with T5 selected).
$ws.Range("A1").Select()
